$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Replicate formatting for the new block of rows (79-91) by copying the
#    existing "Wiremill" header + 4-tier block (rows 66-78) down below the
#    current last block. This brings along the correct style indices for
#    every cell without us needing to know/guess the exact xf ids.
# ---------------------------------------------------------------------------

# Row 79 acts as a new section header, just like row 66 (only B/F/G/H used).
$ws.Range("B66").Copy()
$ws.Range("B79").PasteSpecial(-4122)
$ws.Range("F66").Copy()
$ws.Range("F79").PasteSpecial(-4122)
$ws.Range("G66").Copy()
$ws.Range("G79").PasteSpecial(-4122)
$ws.Range("H66").Copy()
$ws.Range("H79").PasteSpecial(-4122)

$ws.Range("B79").Value = "Building"
$ws.Range("F79").Value = "Tier"
$ws.Range("G79").Value = "Target In/Out/Profit"
$ws.Range("H79").Value = "Target Productivity"

# Rows 80-91: four 3-row tiers, same layout/style as rows 67-78.
$ws.Range("A67:H78").Copy()
$ws.Range("A80").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Fill in values/formulas. The order in which brand-new text values are
#    assigned controls the order new entries are appended to the shared
#    string table, so we deliberately set them in this sequence.
# ---------------------------------------------------------------------------

# New building name first.
$ws.Range("B80").Value = "Motor Industries"

# Tier labels: Tier1, Tier3, Tier4, Tier2 (matches source ordering).
$ws.Range("F80").Value = "Tier 1: Early, No Auto"
$ws.Range("F86").Value = "Tier 3: Hi-Psi, Electric Drive"
$ws.Range("F89").Value = "Tier 4: Complex, Asslines"
$ws.Range("F83").Value = "Tier 2: Steam, Pneumatic"

# Re-use the "Motor Industries" label for the other three tiers.
$ws.Range("B83").Value = "Motor Industries"
$ws.Range("B86").Value = "Motor Industries"
$ws.Range("B89").Value = "Motor Industries"

# --- Block 1: rows 80-82 (Tier 1) ------------------------------------------
$ws.Range("A80").Value = "Base Input"
$ws.Range("C80").Value = 1000
$ws.Range("D80").Value = 5000
$ws.Range("E80").Value = "Empl"
$ws.Range("G80").Value = 600
$ws.Range("H80").Value = 5000

$ws.Range("A81").Value = "Base Output"
$ws.Range("C81").Value = 1620
$ws.Range("G81").Value = 1200

$ws.Range("A82").Value = "Profit"
$ws.Range("C82").Formula = "=C81-C80"
$ws.Range("D82").Formula = "=C82/D80"
$ws.Range("E82").Value = "Per Empl"
$ws.Range("G82").Formula = "=G81-G80"
$ws.Range("H82").Value = "0,12-0,14"

# --- Block 2: rows 83-85 (Tier 2) ------------------------------------------
$ws.Range("A83").Value = "Base Input"
$ws.Range("C83").Formula = "=1600+210"
$ws.Range("D83").Value = 4000
$ws.Range("E83").Value = "Empl"
$ws.Range("G83").Value = 1600
$ws.Range("H83").Value = 4000

$ws.Range("A84").Value = "Base Output"
$ws.Range("C84").Value = 2700
$ws.Range("G84").Value = 2400

$ws.Range("A85").Value = "Profit"
$ws.Range("C85").Formula = "=C84-C83"
$ws.Range("D85").Formula = "=C85/D83"
$ws.Range("E85").Value = "Per Empl"
$ws.Range("G85").Formula = "=G84-G83"
$ws.Range("H85").Value = "0,2-0.23"

# --- Block 3: rows 86-88 (Tier 3) ------------------------------------------
$ws.Range("A86").Value = "Base Input"
$ws.Range("C86").Formula = "=1950+300"
$ws.Range("D86").Value = 3000
$ws.Range("E86").Value = "Empl"
$ws.Range("G86").Value = 2100
$ws.Range("H86").Value = 3000

$ws.Range("A87").Value = "Base Output"
$ws.Range("C87").Value = 3300
$ws.Range("G87").Value = 3200

$ws.Range("A88").Value = "Profit"
$ws.Range("C88").Formula = "=C87-C86"
$ws.Range("D88").Formula = "=C88/D86"
$ws.Range("E88").Value = "Per Empl"
$ws.Range("G88").Formula = "=G87-G86"
$ws.Range("H88").Value = "0,36-0.4"

# --- Block 4: rows 89-91 (Tier 4) ------------------------------------------
$ws.Range("A89").Value = "Base Input"
$ws.Range("C89").Formula = "=2290+495"
$ws.Range("D89").Value = 2000
$ws.Range("E89").Value = "Empl"
$ws.Range("G89").Value = 2400
$ws.Range("H89").Value = 2000

$ws.Range("A90").Value = "Base Output"
$ws.Range("C90").Value = 4620
$ws.Range("G90").Value = 4200

$ws.Range("A91").Value = "Profit"
$ws.Range("C91").Formula = "=C90-C89"
$ws.Range("D91").Formula = "=C91/D89"
$ws.Range("E91").Value = "Per Empl"
$ws.Range("G91").Formula = "=G90-G89"
$ws.Range("H91").Value = "0,9-1.0"

# --- Extra helper cell below the table --------------------------------------
$ws.Range("D93").Formula = "=C90/60"

# ---------------------------------------------------------------------------
# 3) Merge the label cells for each new tier block (B and F columns).
# ---------------------------------------------------------------------------
$ws.Range("B80:B82").Merge()
$ws.Range("F80:F82").Merge()
$ws.Range("B83:B85").Merge()
$ws.Range("F83:F85").Merge()
$ws.Range("B86:B88").Merge()
$ws.Range("F86:F88").Merge()
$ws.Range("B89:B91").Merge()
$ws.Range("F89:F91").Merge()

# ---------------------------------------------------------------------------
# 4) Update the view/selection to match the edited state.
# ---------------------------------------------------------------------------
$ws.Range("J82").Select()
